$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (header / condition numbers) - columns B:E take on the values
# that previously lived in columns O, R, AN, AQ (other columns were
# removed from the dataset, these four remain as the "passive" set).
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 ("CON")
$ws.Range("B2").Value = 15.349999989999999
$ws.Range("C2").Value = 7.2999999899999999
$ws.Range("D2").Value = 13.249999989999999
$ws.Range("E2").Value = 16.349999990000001

# Row 3 ("STR")
$ws.Range("B3").Value = 6.8499999899999997
$ws.Range("C3").Value = 17.04999999
$ws.Range("D3").Value = 23.79999999
$ws.Range("E3").Value = 13.849999989999999

# Update the selection to match the new, narrower range of interest.
$ws.Activate()
$ws.Range("B1:E3").Select() | Out-Null
